$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8395434617996216
$ws.Range("B1").Value = 1.381192803382874
$ws.Range("C1").Value = 1.669219136238098
$ws.Range("D1").Value = 4.820402145385742
$ws.Range("E1").Value = 4.851388454437256
